$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 124.25
$ws.Range("J33").Value = 122.333336
$ws.Range("L33").Value = 122.333336
$ws.Range("N33").Value = -580.333336

$ws.Range("H40").Value = 1013.4894
$ws.Range("I40").Value = 994.0244
$ws.Range("K40").Value = 994.0244
$ws.Range("M40").Value = -819.0244

$ws.Range("H112").Value = 1549.3077
$ws.Range("J112").Value = 1567.3055
$ws.Range("L112").Value = 4701.916499999999
$ws.Range("N112").Value = -6917.916499999999

$ws.Range("H138").Value = 1544.6184
$ws.Range("I138").Value = 1229.7213
$ws.Range("J138").Value = 2825.2
$ws.Range("K138").Value = 3689.1639
$ws.Range("L138").Value = 8475.599999999999
$ws.Range("M138").Value = 1450.8361
$ws.Range("N138").Value = -18755.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3516.4285
$ws.Range("I32").Value = 2898.597
$ws.Range("K32").Value = 2898.597
$ws.Range("M32").Value = -2611.597

$ws.Range("H45").Value = 1689.6086
$ws.Range("I45").Value = 1556.25
$ws.Range("K45").Value = 1556.25
$ws.Range("M45").Value = -1179.25

$ws.Range("H61").Value = 1665.75
$ws.Range("I61").Value = 784.6667
$ws.Range("J61").Value = 7833.3335
$ws.Range("K61").Value = 784.6667
$ws.Range("L61").Value = 7833.3335
$ws.Range("M61").Value = -572.6667
$ws.Range("N61").Value = -8257.333500000001

$ws.Range("H74").Value = 1222.5227
$ws.Range("I74").Value = 909.9091
$ws.Range("J74").Value = 2160.3635
$ws.Range("K74").Value = 909.9091
$ws.Range("L74").Value = 2160.3635
$ws.Range("M74").Value = -35.90909999999997
$ws.Range("N74").Value = -3908.3635

$ws.Range("H77").Value = 1222.5227
$ws.Range("I77").Value = 909.9091
$ws.Range("J77").Value = 2160.3635
$ws.Range("K77").Value = 4549.5455
$ws.Range("L77").Value = 10801.8175
$ws.Range("M77").Value = -181.5455000000002
$ws.Range("N77").Value = -19537.8175

$ws.Range("H132").Value = 1227.3518
$ws.Range("I132").Value = 1021.48834
$ws.Range("J132").Value = 2032.091
$ws.Range("K132").Value = 3064.46502
$ws.Range("L132").Value = 6096.272999999999
$ws.Range("M132").Value = -534.4650200000001
$ws.Range("N132").Value = -11156.273

$ws.Range("H136").Value = 1665.75
$ws.Range("I136").Value = 784.6667
$ws.Range("J136").Value = 7833.3335
$ws.Range("K136").Value = 2354.0001
$ws.Range("L136").Value = 23500.0005
$ws.Range("M136").Value = 195.9998999999998
$ws.Range("N136").Value = -28600.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 402870.34
$ws.Range("I86").Value = 591477.0600000001
$ws.Range("J86").Value = 156230.77
$ws.Range("K86").Value = 591477.0600000001
$ws.Range("L86").Value = 156230.77
$ws.Range("M86").Value = -590354.0600000001
$ws.Range("N86").Value = -158476.77

$ws.Range("H89").Value = 402870.34
$ws.Range("I89").Value = 591477.0600000001
$ws.Range("J89").Value = 156230.77
$ws.Range("K89").Value = 2957385.3
$ws.Range("L89").Value = 781153.85
$ws.Range("M89").Value = -2951769.3
$ws.Range("N89").Value = -792385.85

$ws.Range("H105").Value = 2297.4138
$ws.Range("I105").Value = 2275
$ws.Range("K105").Value = 2275
$ws.Range("M105").Value = -528

$ws.Range("H134").Value = 9719
$ws.Range("I134").Value = 10134.412
$ws.Range("K134").Value = 30403.236
$ws.Range("M134").Value = -27868.236

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1741359.9
$ws.Range("I58").Value = 2899259.2
$ws.Range("J58").Value = 4510.9
$ws.Range("K58").Value = 2899259.2
$ws.Range("L58").Value = 4510.9
$ws.Range("M58").Value = -2899056.2
$ws.Range("N58").Value = -4916.9

$ws.Range("H134").Value = 1953.2285
$ws.Range("I134").Value = 1751.8572
$ws.Range("J134").Value = 2758.7144
$ws.Range("K134").Value = 5255.571599999999
$ws.Range("L134").Value = 8276.143199999999
$ws.Range("M134").Value = -2720.571599999999
$ws.Range("N134").Value = -13346.1432

$ws.Range("H136").Value = 1741359.9
$ws.Range("I136").Value = 2899259.2
$ws.Range("J136").Value = 4510.9
$ws.Range("K136").Value = 8697777.600000001
$ws.Range("L136").Value = 13532.7
$ws.Range("M136").Value = -8695227.600000001
$ws.Range("N136").Value = -18632.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 11744.203
$ws.Range("J131").Value = 15220.796
$ws.Range("L131").Value = 45662.388
$ws.Range("N131").Value = -55742.388

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2246.375
$ws.Range("I102").Value = 2229.4666
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 2229.4666
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -607.4666000000002
$ws.Range("N102").Value = -5744

$ws.Range("H126").Value = 2177801.8
$ws.Range("I126").Value = 7939615.5
$ws.Range("J126").Value = 55028.316
$ws.Range("K126").Value = 23818846.5
$ws.Range("L126").Value = 165084.948
$ws.Range("M126").Value = -23816376.5
$ws.Range("N126").Value = -170024.948

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 986.8421
$ws.Range("J93").Value = 1385.25
$ws.Range("L93").Value = 1385.25
$ws.Range("N93").Value = -3881.25

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 33032.715
$ws.Range("J135").Value = 33032.715
$ws.Range("L135").Value = 33032.715
$ws.Range("N135").Value = -43172.715

$ws.Range("H136").Value = 2305.762
$ws.Range("I136").Value = 1308.3125
$ws.Range("K136").Value = 3924.9375
$ws.Range("M136").Value = -1374.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5149.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5149.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5149.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6397.5

$ws.Range("H65").Value = 5149.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5149.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25747.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31987.5

$ws.Range("H107").Value = 537.7143
$ws.Range("J107").Value = 1088.125
$ws.Range("L107").Value = 3264.375
$ws.Range("N107").Value = -7104.375

$ws.Range("H113").Value = 581.8333
$ws.Range("I113").Value = 283.85715
$ws.Range("K113").Value = 851.5714499999999
$ws.Range("M113").Value = 1318.42855
